# Apply the "Updated BE and ZG" change to the covid19 cases / fatalities workbook.
$wb = $excel.ActiveWorkbook

$wsCases      = $wb.Worksheets.Item(1)   # covid19_cases_switzerland
$wsFatalities = $wb.Worksheets.Item(2)   # covid19_fatalities_switzerland

# ---------------------------------------------------------------------------
# Sheet 1: covid19_cases_switzerland
# ---------------------------------------------------------------------------

# A previously missing SO (Solothurn) figure for 2020-03-24 (row 20) is filled in,
# and the CH total formula no longer needs to add S19 separately since S20 is
# now populated directly.
$wsCases.Range("S20").Value = 129
$wsCases.Range("AB20").Formula = "=SUM(B20:AA20)+X19"

# New row for 2020-03-25 with the updated BE and ZG case counts.
$wsCases.Range("A21").Value = 43915
$wsCases.Range("A21").NumberFormat = "yyyy\-mm\-dd;@"
$wsCases.Range("E21").Value = 624
$wsCases.Range("Z21").Value = 80
$wsCases.Range("AB21").Formula = "=AB20+E21-E20+Z21-Z20"
$wsCases.Range("AB21").NumberFormat = "0"

# ---------------------------------------------------------------------------
# Sheet 2: covid19_fatalities_switzerland
# ---------------------------------------------------------------------------

# New row for 2020-03-25 with the updated fatalities per canton.
$wsFatalities.Range("A21").Value = 43915
$wsFatalities.Range("A21").NumberFormat = "yyyy\-mm\-dd;@"
$wsFatalities.Range("B21").Value = 2
$wsFatalities.Range("D21").Value = 2
$wsFatalities.Range("E21").Value = 6
$wsFatalities.Range("F21").Value = 4
$wsFatalities.Range("G21").Value = 5
$wsFatalities.Range("H21").Value = 1
$wsFatalities.Range("I21").Value = 12
$wsFatalities.Range("K21").Value = 6
$wsFatalities.Range("N21").Value = 2
$wsFatalities.Range("Q21").Value = 1
$wsFatalities.Range("S21").Value = 1
$wsFatalities.Range("U21").Value = 1
$wsFatalities.Range("V21").Value = 53
$wsFatalities.Range("X21").Value = 17
$wsFatalities.Range("Y21").Value = 13
$wsFatalities.Range("AA21").Value = 5
$wsFatalities.Range("AB21").Formula = "=SUM(B21:AA21)"
$wsFatalities.Range("AB21").NumberFormat = "0"

# ---------------------------------------------------------------------------
# Recalculate so cached formula results are up to date.
# ---------------------------------------------------------------------------
$excel.Calculate()

# ---------------------------------------------------------------------------
# Restore sheet/selection state: sheet1 becomes the active tab again with a
# selection on M25, while sheet2's selection moves to A28.
# ---------------------------------------------------------------------------
$wsFatalities.Activate()
[void]$wsFatalities.Range("A28").Select()
$wsCases.Activate()
[void]$wsCases.Range("M25").Select()
